$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.018146673569414
$ws.Range("D2").Value = 1.024205015779447
$ws.Range("E2").Value = 1.046848230342142
$ws.Range("F2").Value = 1.050579995130683
$ws.Range("I2").Value = 1.030058161607156
$ws.Range("J2").Value = 1.023356762531119
$ws.Range("K2").Value = 1.027034201238806
$ws.Range("L2").Value = 1.049612517975328
$ws.Range("M2").Value = 1.053333878358224
$ws.Range("N2").Value = 1.011896732055202
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.01895491357453
$ws.Range("D3").Value = 1.024781769266521
$ws.Range("E3").Value = 1.047886765717097
$ws.Range("F3").Value = 1.051632725503216
$ws.Range("I3").Value = 1.030192409588552
$ws.Range("J3").Value = 1.023801938309647
$ws.Range("K3").Value = 1.027418623519208
$ws.Range("L3").Value = 1.050462189824797
$ws.Range("M3").Value = 1.054198464564233
$ws.Range("N3").Value = 1.012043665117992
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.019478115695039
$ws.Range("D4").Value = 1.025154770441733
$ws.Range("E4").Value = 1.048559736440967
$ws.Range("F4").Value = 1.052314578351837
$ws.Range("I4").Value = 1.030277488851436
$ws.Range("J4").Value = 1.024089570743041
$ws.Range("K4").Value = 1.0276664745932
$ws.Range("L4").Value = 1.051012348357931
$ws.Range("M4").Value = 1.0547579734125
$ws.Range("N4").Value = 1.012138588311127
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.019698120345175
$ws.Range("D5").Value = 1.025311531306125
$ws.Range("E5").Value = 1.048842884398875
$ws.Range("F5").Value = 1.05260138770807
$ws.Range("I5").Value = 1.030312827277336
$ws.Range("J5").Value = 1.024210388337002
$ws.Range("K5").Value = 1.027770455508196
$ws.Range("L5").Value = 1.051243721338312
$ws.Range("M5").Value = 1.054993205230426
$ws.Range("N5").Value = 1.012178457176797
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.01973506301272
$ws.Range("D6").Value = 1.025337849232278
$ws.Range("E6").Value = 1.048890439710933
$ws.Range("F6").Value = 1.052649553521719
$ws.Range("I6").Value = 1.030318735580296
$ws.Range("J6").Value = 1.024230668066433
$ws.Range("K6").Value = 1.02778790167872
$ws.Range("L6").Value = 1.05128257491872
$ws.Range("M6").Value = 1.055032702498106
$ws.Range("N6").Value = 1.012185149163158
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.019481055211382
$ws.Range("D7").Value = 1.025156865282085
$ws.Range("E7").Value = 1.048563518971021
$ws.Range("F7").Value = 1.05231841009127
$ws.Range("I7").Value = 1.030277962731502
$ws.Range("J7").Value = 1.024091185520046
$ws.Range("K7").Value = 1.027667864840266
$ws.Range("L7").Value = 1.051015439636812
$ws.Range("M7").Value = 1.054761116536597
$ws.Range("N7").Value = 1.012139121186113
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.018419776321182
$ws.Range("D8").Value = 1.024399972074914
$ws.Range("E8").Value = 1.047199007084129
$ws.Range("F8").Value = 1.050935631807771
$ws.Range("I8").Value = 1.03010390110678
$ws.Range("J8").Value = 1.023507299340611
$ws.Range("K8").Value = 1.027164303314424
$ws.Range("L8").Value = 1.049899592688768
$ws.Range("M8").Value = 1.053626055531801
$ws.Range("N8").Value = 1.011946420084444
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.016551390287422
$ws.Range("D9").Value = 1.02306479592698
$ws.Range("E9").Value = 1.044802022460878
$ws.Range("F9").Value = 1.048504143781331
$ws.Range("I9").Value = 1.029783522555148
$ws.Range("J9").Value = 1.022475201661363
$ws.Range("K9").Value = 1.026270155546002
$ws.Range("L9").Value = 1.047936151749639
$ws.Range("M9").Value = 1.051626463805487
$ws.Range("N9").Value = 1.011605706130513
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.015307041464747
$ws.Range("D10").Value = 1.022173814345917
$ws.Range("E10").Value = 1.043209099243703
$ws.Range("F10").Value = 1.046886667783982
$ws.Range("I10").Value = 1.029560800040062
$ws.Range("J10").Value = 1.02178503796283
$ws.Range("K10").Value = 1.025669550662573
$ws.Range("L10").Value = 1.046629136759671
$ws.Range("M10").Value = 1.050293815253129
$ws.Range("N10").Value = 1.01137781252228
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.014768536535882
$ws.Range("D11").Value = 1.021787825390617
$ws.Range("E11").Value = 1.0425205578577
$ws.Range("F11").Value = 1.046187128328855
$ws.Range("I11").Value = 1.029462201374751
$ws.Range("J11").Value = 1.021485704888454
$ws.Range("K11").Value = 1.025408428033314
$ws.Range("L11").Value = 1.046063655754016
$ws.Range("M11").Value = 1.049716872837426
$ws.Range("N11").Value = 1.011278958577203
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.01456855920891
$ws.Range("D12").Value = 1.021644425320105
$ws.Range("E12").Value = 1.042264984715457
$ws.Range("F12").Value = 1.045927415042789
$ws.Range("I12").Value = 1.029425253922752
$ws.Range("J12").Value = 1.021374447031106
$ws.Range("K12").Value = 1.02531127785842
$ws.Range("L12").Value = 1.045853681511796
$ws.Range("M12").Value = 1.049502587086027
$ws.Range("N12").Value = 1.011242213967003
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.014611452858518
$ws.Range("D13").Value = 1.021675186298441
$ws.Range("E13").Value = 1.042319797783507
$ws.Range("F13").Value = 1.045983118671724
$ws.Range("I13").Value = 1.029433193904486
$ws.Range("J13").Value = 1.021398315470947
$ws.Range("K13").Value = 1.025332124014048
$ws.Range("L13").Value = 1.045898718494252
$ws.Range("M13").Value = 1.04954855135365
$ws.Range("N13").Value = 1.011250096974714
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.01475200537325
$ws.Range("D14").Value = 1.021775972428804
$ws.Range("E14").Value = 1.042499428389218
$ws.Range("F14").Value = 1.04616565775562
$ws.Range("I14").Value = 1.029459153882183
$ws.Range("J14").Value = 1.021476509746662
$ws.Range("K14").Value = 1.02540040077793
$ws.Range("L14").Value = 1.046046297773737
$ws.Range("M14").Value = 1.049699159571938
$ws.Range("N14").Value = 1.011275921781563
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.014838610735394
$ws.Range("D15").Value = 1.021838066616821
$ws.Range("E15").Value = 1.042610128845313
$ws.Range("F15").Value = 1.046278142939844
$ws.Range("I15").Value = 1.029475105834361
$ws.Range("J15").Value = 1.021524678280992
$ws.Range("K15").Value = 1.025442447496154
$ws.Range("L15").Value = 1.046137235612358
$ws.Range("M15").Value = 1.049791956457782
$ws.Range("N15").Value = 1.011291829882421
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.015342786833271
$ws.Range("D16").Value = 1.022199427342061
$ws.Range("E16").Value = 1.043254820918604
$ws.Range("F16").Value = 1.046933111650142
$ws.Range("I16").Value = 1.029567298298943
$ws.Range("J16").Value = 1.021804893527724
$ws.Range("K16").Value = 1.025686858346195
$ws.Range("L16").Value = 1.04666667573101
$ws.Range("M16").Value = 1.050332107267611
$ws.Range("N16").Value = 1.011384369487489
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.015659126150684
$ws.Range("D17").Value = 1.022426050210369
$ws.Range("E17").Value = 1.043659542344217
$ws.Range("F17").Value = 1.047344181272255
$ws.Range("I17").Value = 1.02962455087606
$ws.Range("J17").Value = 1.021980535327085
$ws.Range("K17").Value = 1.025839888714376
$ws.Range("L17").Value = 1.046998904549283
$ws.Range("M17").Value = 1.050670957826529
$ws.Range("N17").Value = 1.011442370673676
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.015843670987197
$ws.Range("D18").Value = 1.022558217239845
$ws.Range("E18").Value = 1.043895725706223
$ws.Range("F18").Value = 1.04758403205619
$ws.Range("I18").Value = 1.029657737036471
$ws.Range("J18").Value = 1.022082937097705
$ws.Range("K18").Value = 1.025929046700365
$ws.Range("L18").Value = 1.047192732827177
$ws.Range("M18").Value = 1.050868613525362
$ws.Range("N18").Value = 1.011476184922348
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.015906600963727
$ws.Range("D19").Value = 1.02260327963413
$ws.Range("E19").Value = 1.043976277815702
$ws.Range("F19").Value = 1.047665828658751
$ws.Range("I19").Value = 1.029669017305543
$ws.Range("J19").Value = 1.022117845431339
$ws.Range("K19").Value = 1.025959429926591
$ws.Range("L19").Value = 1.047258830843372
$ws.Range("M19").Value = 1.050936010664418
$ws.Range("N19").Value = 1.01148771183612
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.015625182878212
$ws.Range("D20").Value = 1.022401737597526
$ws.Range("E20").Value = 1.043616107540601
$ws.Range("F20").Value = 1.047300069025376
$ws.Range("I20").Value = 1.029618429756294
$ws.Range("J20").Value = 1.021961695490534
$ws.Range("K20").Value = 1.025823480550065
$ws.Range("L20").Value = 1.046963254899246
$ws.Range("M20").Value = 1.050634601338163
$ws.Range("N20").Value = 1.011436149436785
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.014710614842002
$ws.Range("D21").Value = 1.021746294146017
$ws.Range("E21").Value = 1.042446526628664
$ws.Range("F21").Value = 1.046111901037416
$ws.Range("I21").Value = 1.029451518240456
$ws.Range("J21").Value = 1.021453485466885
$ws.Range("K21").Value = 1.025380299314065
$ws.Range("L21").Value = 1.046002837400193
$ws.Range("M21").Value = 1.049654808747325
$ws.Range("N21").Value = 1.011268317730377
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.014135864419119
$ws.Range("D22").Value = 1.021334037831645
$ws.Range("E22").Value = 1.041712216826956
$ws.Range("F22").Value = 1.045365586805021
$ws.Range("I22").Value = 1.029344702940957
$ws.Range("J22").Value = 1.021133536288737
$ws.Range("K22").Value = 1.025100742404082
$ws.Range("L22").Value = 1.045399393511427
$ws.Range("M22").Value = 1.049038869273865
$ws.Range("N22").Value = 1.011162645874159
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.014440523991724
$ws.Range("D23").Value = 1.021552596652772
$ws.Range("E23").Value = 1.042101388359292
$ws.Range("F23").Value = 1.045761152222855
$ws.Range("I23").Value = 1.029401504890523
$ws.Range("J23").Value = 1.021303186618078
$ws.Range("K23").Value = 1.025249026803866
$ws.Range("L23").Value = 1.045719251557709
$ws.Range("M23").Value = 1.049365381153747
$ws.Range("N23").Value = 1.011218678564715
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.015640520291061
$ws.Range("D24").Value = 1.022412723477341
$ws.Range("E24").Value = 1.043635733499691
$ws.Range("F24").Value = 1.047320001201219
$ws.Range("I24").Value = 1.029621196270513
$ws.Range("J24").Value = 1.021970208547415
$ws.Range("K24").Value = 1.025830895008639
$ws.Range("L24").Value = 1.046979363303071
$ws.Range("M24").Value = 1.050651029239738
$ws.Range("N24").Value = 1.011438960598084
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.017034200074414
$ws.Range("D25").Value = 1.023410130546142
$ws.Range("E25").Value = 1.045420811177712
$ws.Range("F25").Value = 1.049132126209569
$ws.Range("I25").Value = 1.029867961956948
$ws.Range("J25").Value = 1.022742398167697
$ws.Range("K25").Value = 1.02650211362723
$ws.Range("L25").Value = 1.048443408808288
$ws.Range("M25").Value = 1.052143337554652
$ws.Range("N25").Value = 1.011693922839263
